$d = $word.ActiveDocument
$CR = [char]13

# ---------------------------------------------------------------------------
# Change 1: split the "About the city of Toronto" heading paragraph into two
# paragraphs - a new "Introduction" heading (bigger, 18pt) followed by the
# existing "About the city of Toronto" heading shrunk down to 14pt.
# ---------------------------------------------------------------------------
$aboutIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd($CR) -eq "About the city of Toronto") {
        $aboutIdx = $i
        break
    }
}

$d.Paragraphs.Item($aboutIdx).Range.InsertParagraphBefore()

# after the insert, the new blank paragraph takes the old index and the
# "About the city of Toronto" paragraph is pushed one slot later
$introPara = $d.Paragraphs.Item($aboutIdx)
$aboutPara = $d.Paragraphs.Item($aboutIdx + 1)

$introPara.Range.Text = "Introduction"
$introPara.Range.Font.Size = 18
$introPara.Range.Font.SizeBi = 18
$introPara.Range.LanguageID = "en-US"

$aboutPara.Range.Font.Size = 14
$aboutPara.Range.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# Change 2: rewrite the "Now we want to concentrate..." paragraph as a
# "Problem:-" lead-in followed by the reworded sentence.
# ---------------------------------------------------------------------------
$target = "Now we want to concentrate on the business problem that we want to solve and that is to find the most optimum location or neighborhood to set up a new Indian Restaurant in Toronto. Is it possible to predict the success of a new restaurant even before establishing it? Who would benefit from this Project?"
$problemPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd($CR) -eq $target) {
        $problemPara = $cand
        break
    }
}

$problemPara.Range.Text = "Problem:- The business problem that we want to solve is that to find the most optimum location or neighborhood to set up a new Indian Restaurant in Toronto. Is it possible to predict the success of a new restaurant even before establishing it? Who would benefit from this Project?"

$paraStart = $problemPara.Range.Start

$problemRange = $d.Range($paraStart, $paraStart + 7)
$problemRange.Bold = 1
$problemRange.LanguageID = "en-US"

$colonRange = $d.Range($paraStart + 7, $paraStart + 9)
$colonRange.LanguageID = "en-US"

$spaceTRange = $d.Range($paraStart + 9, $paraStart + 11)
$spaceTRange.LanguageID = "en-US"

$isRange = $d.Range($paraStart + 53, $paraStart + 55)
$isRange.LanguageID = "en-US"

Write-Output "done"
